# Daily attendance processing - 2025-12-23 21:30:53
#
# For every data row in column "G" ("Recorded By") whose text begins with
# the literal prefix "System, ", move that leading "System" token to the
# end of the comma-separated list (e.g. "System, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$prefix = "System, "

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $value = $cell.Value()

    if ($value -ne $null -and $value.StartsWith($prefix)) {
        $rest = $value.Substring($prefix.Length)
        $cell.Value = $rest + ", System"
    }
}
